$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.944.51'
$ws.Range('E2').Value = '  -1.07%  '

$ws.Range('D3').Value = '1.638.10'
$ws.Range('E3').Value = '  -0.43%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.41'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.84%  '

$ws.Range('E6').Value = '  -0.26%  '

$ws.Range('E7').Value = '  +0.33%  '

$ws.Range('E8').Value = '  -0.74%  '

$ws.Range('E9').Value = '  +0.18%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.62'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.91%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0794'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.15%  '

$ws.Range('D12').Value = '1.864.90'
$ws.Range('E12').Value = '  -0.44%  '

$ws.Range('E13').Value = '  -1.14%  '

$ws.Range('D14').Value = '1.638.49'
$ws.Range('E14').Value = '  -0.35%  '

$ws.Range('E15').Value = '  -0.93%  '

$ws.Range('D16').Value = '0.0₃0764'
$ws.Range('E16').Value = '  -0.13%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '62.95'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.93%  '

$ws.Range('D18').Value = '25.964.66'
$ws.Range('E18').Value = '  -0.94%  '

$ws.Range('E19').Value = '  +0.26%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '192.94'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.31%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.92'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.38%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.28'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.84%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '144.18'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.46%  '

$ws.Range('E25').Value = '  +1.16%  '

$ws.Range('E26').Value = '  +0.28%  '

$ws.Range('E27').Value = '  +3.62%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.85'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.37%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.55'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.37%  '

$ws.Range('E30').Value = '  -0.72%  '

$ws.Range('E31').Value = '  -0.22%  '

$ws.Range('E32').Value = '  -1.03%  '

$ws.Range('E34').Value = '  -4.25%  '

$ws.Range('E35').Value = '  +1.75%  '

$ws.Range('E36').Value = '  -1.22%  '

$ws.Range('D37').Value = '1.138.97'
$ws.Range('E37').Value = '  +0.42%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.545'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.64%  '

$ws.Range('E39').Value = '  -1.06%  '

$ws.Range('E40').Value = '  +0.38%  '

$ws.Range('E41').Value = '  +0.28%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.48'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.08%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '99.29'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.80%  '

$ws.Range('D45').Value = '1.774.54'
$ws.Range('E45').Value = '  -0.45%  '

$ws.Range('E46').Value = '  +2.37%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '56.63'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.66%  '

$ws.Range('E48').Value = '  +3.16%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.48'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.21%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.65'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.75%  '

$ws.Range('E51').Value = '  -0.72%  '
